$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Step 1: write all new cell VALUES first (values before formatting avoids a
#     COUNTIF/whole-column dependency staleness bug triggered by pasting formats
#     onto a cell before its value is set) ---

# Row 107
$ws.Cells.Item(107,1).Value = 'The Misuse of BLUP in Ecology and Evolution'
$ws.Cells.Item(107,2).Value = 'The american naturalist'
$ws.Cells.Item(107,3).Value = 2010
$ws.Cells.Item(107,4).Value = 'Hadfield'
$ws.Cells.Item(107,6).Value = 'yes'
$ws.Cells.Item(107,7).Value = 'BLUPs'
$ws.Cells.Item(107,9).Value = 'yes'
$ws.Cells.Item(107,10).Value = 'yes'
$ws.Cells.Item(107,11).Value = 'no'

# Row 108
$ws.Cells.Item(108,1).Value = 'Measuring growth patterns in the field: effects of sampling regime and methods on standardized estimates'
$ws.Cells.Item(108,2).Value = 'Canadian Journal of Zoology'
$ws.Cells.Item(108,3).Value = 2011
$ws.Cells.Item(108,4).Value = 'Martin & Pelletier'
$ws.Cells.Item(108,6).Value = 'yes'
$ws.Cells.Item(108,7).Value = 'BLUPs'
$ws.Cells.Item(108,9).Value = 'yes'
$ws.Cells.Item(108,10).Value = 'yes'
$ws.Cells.Item(108,11).Value = 'no'

# Row 109
$ws.Cells.Item(109,1).Value = 'Organisms as Ecosystem Engineers'
$ws.Cells.Item(109,2).Value = 'Oikos'
$ws.Cells.Item(109,3).Value = 1994
$ws.Cells.Item(109,4).Value = 'Jones, Lawton & Shachak'
$ws.Cells.Item(109,6).Value = '≈'
$ws.Cells.Item(109,7).Value = 'Concept of engineers species and example of Beavers'
$ws.Cells.Item(109,9).Value = 'yes'
$ws.Cells.Item(109,10).Value = 'yes'
$ws.Cells.Item(109,11).Value = 'yes'

# Row 110
$ws.Cells.Item(110,1).Value = 'The measurement of selection on correlated characters'
$ws.Cells.Item(110,2).Value = 'Evolution'
$ws.Cells.Item(110,3).Value = 1983
$ws.Cells.Item(110,4).Value = 'Lande & Arnold'
$ws.Cells.Item(110,9).Value = 'yes'
$ws.Cells.Item(110,10).Value = 'yes'
$ws.Cells.Item(110,11).Value = 'no'

# Row 111
$ws.Cells.Item(111,1).Value = 'Behavioral types as predictors of survival in Trinidadian guppies (Poecilia reticulata)'
$ws.Cells.Item(111,2).Value = 'Behavioral Ecology'
$ws.Cells.Item(111,3).Value = 2010
$ws.Cells.Item(111,4).Value = 'Smith et Blumstein'
$ws.Cells.Item(111,6).Value = '≈'
$ws.Cells.Item(111,7).Value = 'bold, active and exploratory guppies survive longer when exposed to predators'
$ws.Cells.Item(111,9).Value = 'yes'
$ws.Cells.Item(111,10).Value = 'yes'
$ws.Cells.Item(111,11).Value = 'yes'

# Row 112
$ws.Cells.Item(112,1).Value = 'Impact of Body Reserves on Energy Expenditure, Water Flux, and Mating Success in Breeding Male Northern Elephant Seals'
$ws.Cells.Item(112,2).Value = 'Physiological and Biochemical Zoology'
$ws.Cells.Item(112,3).Value = 2012
$ws.Cells.Item(112,4).Value = 'Crocker, Houser & Webb'
$ws.Cells.Item(112,6).Value = '≈'
$ws.Cells.Item(112,7).Value = 'Bigger elephant seals have a better mating success'
$ws.Cells.Item(112,9).Value = 'yes'
$ws.Cells.Item(112,10).Value = 'yes'
$ws.Cells.Item(112,11).Value = 'yes'

# Row 113
$ws.Cells.Item(113,1).Value = 'Reproductive success and failure: the role of winter body mass in reproductive allocation in Norwegian moose'
$ws.Cells.Item(113,2).Value = 'Oecologia'
$ws.Cells.Item(113,3).Value = 2013
$ws.Cells.Item(113,4).Value = 'Milner et al.'
$ws.Cells.Item(113,6).Value = '≈'
$ws.Cells.Item(113,7).Value = 'Link between winter body mass variation and reproductive success in female northern moose (Alces alces)'
$ws.Cells.Item(113,9).Value = 'yes'
$ws.Cells.Item(113,10).Value = 'yes'
$ws.Cells.Item(113,11).Value = 'yes'

# K25: mark this citation as not-cited (was 'yes', now 'no'); style turns red like K4
$ws.Cells.Item(25,11).Value = 'no'

# --- Step 2: copy cell FORMATTING from matching template rows (format-only paste) ---

$ws.Range("A103:D103").Copy()
$ws.Range("A107:D107").PasteSpecial(-4122)
$ws.Range("F103:G103").Copy()
$ws.Range("F107:G107").PasteSpecial(-4122)
$ws.Range("I103:K103").Copy()
$ws.Range("I107:K107").PasteSpecial(-4122)
$ws.Rows.Item(107).RowHeight = 17

$ws.Range("A103:D103").Copy()
$ws.Range("A108:D108").PasteSpecial(-4122)
$ws.Range("F103:G103").Copy()
$ws.Range("F108:G108").PasteSpecial(-4122)
$ws.Range("I103:K103").Copy()
$ws.Range("I108:K108").PasteSpecial(-4122)
$ws.Rows.Item(108).RowHeight = 17

$ws.Range("A105:D105").Copy()
$ws.Range("A109:D109").PasteSpecial(-4122)
$ws.Range("F105:G105").Copy()
$ws.Range("F109:G109").PasteSpecial(-4122)
$ws.Range("I105:K105").Copy()
$ws.Range("I109:K109").PasteSpecial(-4122)
$ws.Rows.Item(109).RowHeight = 17

$ws.Range("A106:D106").Copy()
$ws.Range("A110:D110").PasteSpecial(-4122)
$ws.Range("I106:K106").Copy()
$ws.Range("I110:K110").PasteSpecial(-4122)
$ws.Rows.Item(110).RowHeight = 17

$ws.Range("A105:D105").Copy()
$ws.Range("A111:D111").PasteSpecial(-4122)
$ws.Range("F105:G105").Copy()
$ws.Range("F111:G111").PasteSpecial(-4122)
$ws.Range("I105:K105").Copy()
$ws.Range("I111:K111").PasteSpecial(-4122)
$ws.Rows.Item(111).RowHeight = 17

$ws.Range("A105:D105").Copy()
$ws.Range("A112:D112").PasteSpecial(-4122)
$ws.Range("F105:G105").Copy()
$ws.Range("F112:G112").PasteSpecial(-4122)
$ws.Range("I105:K105").Copy()
$ws.Range("I112:K112").PasteSpecial(-4122)
$ws.Rows.Item(112).RowHeight = 17

$ws.Range("A105:D105").Copy()
$ws.Range("A113:D113").PasteSpecial(-4122)
$ws.Range("F105:G105").Copy()
$ws.Range("F113:G113").PasteSpecial(-4122)
$ws.Range("I105:K105").Copy()
$ws.Range("I113:K113").PasteSpecial(-4122)
$ws.Rows.Item(113).RowHeight = 17

# K25 style -> red 'no' look, matching K4's style
$ws.Range("K4:K4").Copy()
$ws.Range("K25:K25").PasteSpecial(-4122)

# --- Step 3: force a full recalculation so cached formula results (N19/N20/N21/N24) are fresh ---
$excel.CalculateFull()

# --- Step 4: update the sheet view's selection / active cell ---
$ws.Range("A113").Select() | Out-Null

